$wb = $excel.ActiveWorkbook

# --- cuota-extra: update selection, leave it as the non-active tab ---
$wsCuotaExtra = $wb.Worksheets.Item("cuota-extra")
$wsCuotaExtra.Activate()
$wsCuotaExtra.Range("G27").Select()

# --- ene2025: record the January payments (65000) for several neighbors ---
$wsEne = $wb.Worksheets.Item("ene2025")
$wsEne.Activate()

$wsEne.Range("C4").Value = 65000
$wsEne.Range("C9").Value = 65000
$wsEne.Range("C10").Value = 65000
$wsEne.Range("D10").Value = 65000
$wsEne.Range("C15").Value = 65000
$wsEne.Range("C16").Value = 65000
$wsEne.Range("C17").Value = 65000
$wsEne.Range("C19").Value = 65000
$wsEne.Range("C20").Value = 65000
$wsEne.Range("C23").Value = 65000

# ene2025 ends up as the active sheet with C23 selected
$wsEne.Range("C23").Select()
